$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of the data rows (A2:C52), leaving the header row (row 1)
# and cell formatting/styles intact. This matches the diff where all values
# in rows 2-52 for columns A, B, C become empty strings.
$ws.Range("A2:C52").ClearContents()
